$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for several rows
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = -2
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = -3
$ws.Range("F10").Value = 3
$ws.Range("F12").Value = -3
$ws.Range("F15").Value = -5
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = 0

# Row 27 updates (E, F, H, I columns)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 8

# More column F updates
$ws.Range("F29").Value = -2
$ws.Range("F32").Value = 0
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 4
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = 2
$ws.Range("F48").Value = 6
$ws.Range("F52").Value = -5
$ws.Range("F57").Value = 2
$ws.Range("F60").Value = -1
